# Update the "Assets" sheet: refresh the account/wealth-class labels that
# feed the Wealth Allocation table, tweak the balances that drive them, and
# drop the now-unused trailing "Fixed Asset" row (its numbers were folded
# into row 6).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Assets")

# --- Financial Account / Asset Type table (A:E) -----------------------
$ws.Range("A3").Value = "AMAR"
$ws.Range("B3").Value = "Liquid"
$ws.Range("C3").Value = 450
$ws.Range("D3").Value = 450

$ws.Range("A4").Value = "MEGA"
$ws.Range("B4").Value = "Fixed"
$ws.Range("C4").Value = 250
$ws.Range("D4").Value = 250

$ws.Range("A5").Value = "LERO"
$ws.Range("B5").Value = "Fixed"
$ws.Range("C5").Value = 225
$ws.Range("D5").Value = 225

$ws.Range("A6").Value = "PLAO"
$ws.Range("B6").Value = "Liquid"
$ws.Range("C6").Value = 350
$ws.Range("D6").Value = 350

# --- Wealth Class table (G:J) ------------------------------------------
$ws.Range("G3").Value = "MARC"
$ws.Range("H3").Value = 300
$ws.Range("I3").Value = 300

$ws.Range("G4").Value = "HATO"
$ws.Range("H4").Value = 350
$ws.Range("I4").Value = 350
$ws.Range("G4:J4").Style = $ws.Range("G5").Style

$ws.Range("G5").Value = "NHA"
$ws.Range("H5").Value = 150
$ws.Range("I5").Value = 150

$ws.Range("G6").Value = "Fixed Asset"
$ws.Range("H6").Value = 475
$ws.Range("I6").Value = 475
$ws.Range("G6:J6").Style = $ws.Range("G5").Style

# Row 7 ("Fixed Asset" / 234 / 234) has been merged into row 6 above, so
# drop it entirely (shifts the table's bottom edge up one row).
$ws.Rows.Item(7).Delete()

# --- Asset Type / Asset Total summary table (L:M) -----------------------
$ws.Range("M2").Value = 800
$ws.Range("L3").Value = "Fixed"
$ws.Range("M3").Value = 475
